$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 2): restyle "Suma" header (H2) and add new "Total" header (I2) ---
# H2 keeps its text (shared string "Suma") but gets a green fill added (bold+border+center+green)
$ws.Range("C2").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$ws.Range("H2").Interior.Color = 5296274

# I2 is a brand-new header cell "Total" with a green fill, regular font, bordered, no special alignment
$ws.Range("C5").Copy()
$ws.Range("I2").PasteSpecial(-4122)
$ws.Range("I2").Interior.Color = 5296274
$ws.Range("I2").Value = "Total"

# --- Row 3: restyle H3 (Suma total) to the bold/green style, and add I3 = 10 (new "Total" input) ---
$ws.Range("B3").Copy()
$ws.Range("H3").PasteSpecial(-4122)

$ws.Range("C3").Copy()
$ws.Range("I3").PasteSpecial(-4122)
$ws.Range("I3").Value = 10

# --- Row 9: switch the "Suma" formulas from referencing H3 to referencing the new Total cell $I$3 ---
$ws.Range("C9").Formula = '=$I$3-C3'
$ws.Range("D9").Formula = '=$I$3-D3'
$ws.Range("E9").Formula = '=$I$3-E3'
$ws.Range("F9").Formula = '=$I$3-F3'
$ws.Range("G9").Formula = '=$I$3-G3'

# Move the active selection, matching where the author ended up working
[void]$ws.Range("C10").Select()
